$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for the data range so numeric-looking values
# (prices, percentages, hour numbers) are stored as text, matching the
# original inline-string cell contents rather than being auto-converted
# to numbers by Excel.
$ws.Range("B2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "292.32"
$ws.Range("E2").Value = "0.56%"
$ws.Range("G2").Value = "18"
$ws.Range("D3").Value = "30.98"
$ws.Range("E3").Value = "0.38%"
$ws.Range("G3").Value = "18"
$ws.Range("D4").Value = "4.953"
$ws.Range("E4").Value = "1.06%"
$ws.Range("G4").Value = "18"
$ws.Range("D5").Value = "0.07473"
$ws.Range("E5").Value = "3.14%"
$ws.Range("G5").Value = "18"
$ws.Range("E6").Value = "-5.83%"
$ws.Range("G6").Value = "18"
$ws.Range("D7").Value = "7.769"
$ws.Range("E7").Value = "1.27%"
$ws.Range("G7").Value = "18"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9212"
$ws.Range("E8").Value = "2.53%"
$ws.Range("G8").Value = "18"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.09347"
$ws.Range("E9").Value = "18.67%"
$ws.Range("G9").Value = "18"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1732"
$ws.Range("E10").Value = "4.10%"
$ws.Range("G10").Value = "18"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08325"
$ws.Range("E11").Value = "2.38%"
$ws.Range("G11").Value = "18"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03263"
$ws.Range("E12").Value = "5.21%"
$ws.Range("G12").Value = "18"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09938"
$ws.Range("E13").Value = "-0.71%"
$ws.Range("G13").Value = "18"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001495"
$ws.Range("E14").Value = "-0.37%"
$ws.Range("G14").Value = "18"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005709"
$ws.Range("E15").Value = "-2.08%"
$ws.Range("G15").Value = "18"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.474"
$ws.Range("E16").Value = "0.31%"
$ws.Range("G16").Value = "18"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.765"
$ws.Range("E17").Value = "1.47%"
$ws.Range("G17").Value = "18"
$ws.Range("D18").Value = "2.181"
$ws.Range("E18").Value = "5.17%"
$ws.Range("G18").Value = "18"
$ws.Range("E19").Value = "0.41%"
$ws.Range("G19").Value = "18"
$ws.Range("D20").Value = "0.1305"
$ws.Range("E20").Value = "0.55%"
$ws.Range("G20").Value = "18"
$ws.Range("D21").Value = "4.138"
$ws.Range("E21").Value = "4.33%"
$ws.Range("G21").Value = "18"
$ws.Range("D22").Value = "0.2118"
$ws.Range("E22").Value = "-8.09%"
$ws.Range("G22").Value = "18"
$ws.Range("D23").Value = "0.04522"
$ws.Range("E23").Value = "0.18%"
$ws.Range("G23").Value = "18"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").Value = "0.53%"
$ws.Range("G24").Value = "18"
$ws.Range("D25").Value = "0.004277"
$ws.Range("E25").Value = "-2.79%"
$ws.Range("G25").Value = "18"
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").Value = "-0.39%"
$ws.Range("G26").Value = "18"
$ws.Range("D27").Value = "0.0003385"
$ws.Range("E27").Value = "-0.37%"
$ws.Range("G27").Value = "18"
$ws.Range("G28").Value = "18"
$ws.Range("G29").Value = "18"
$ws.Range("G30").Value = "18"
$ws.Range("G31").Value = "18"
$ws.Range("G32").Value = "18"
$ws.Range("G33").Value = "18"
$ws.Range("G34").Value = "18"
$ws.Range("G35").Value = "18"
$ws.Range("G36").Value = "18"
$ws.Range("G37").Value = "18"
$ws.Range("G38").Value = "18"
$ws.Range("D39").Value = "0.01625"
$ws.Range("E39").Value = "3.69%"
$ws.Range("G39").Value = "18"
$ws.Range("D40").Value = "0.04574"
$ws.Range("E40").Value = "4.69%"
$ws.Range("G40").Value = "18"
$ws.Range("D41").Value = "0.007440"
$ws.Range("E41").Value = "1.62%"
$ws.Range("G41").Value = "18"
$ws.Range("D42").Value = "0.009814"
$ws.Range("E42").Value = "-1.74%"
$ws.Range("G42").Value = "18"
$ws.Range("D43").Value = "0.1361"
$ws.Range("E43").Value = "3.65%"
$ws.Range("G43").Value = "18"
$ws.Range("D44").Value = "0.002153"
$ws.Range("E44").Value = "6.74%"
$ws.Range("G44").Value = "18"
$ws.Range("D45").Value = "0.01004"
$ws.Range("E45").Value = "5.74%"
$ws.Range("G45").Value = "18"
$ws.Range("D46").Value = "0.00006091"
$ws.Range("E46").Value = "6.26%"
$ws.Range("G46").Value = "18"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.37%"
$ws.Range("G47").Value = "18"
$ws.Range("D48").Value = "2.654"
$ws.Range("E48").Value = "18.44%"
$ws.Range("G48").Value = "18"
$ws.Range("D49").Value = "0.001994"
$ws.Range("E49").Value = "-31.26%"
$ws.Range("G49").Value = "18"
$ws.Range("D50").Value = "0.00002094"
$ws.Range("E50").Value = "-0.37%"
$ws.Range("G50").Value = "18"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").Value = "-0.37%"
$ws.Range("G51").Value = "18"

Write-Output "Applied 146 cell updates"
